# Update "想去人数" (want-to-go count) figures across all sheets to the
# freshly scraped values (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2621
$ws1.Range("F3").Value = 568
$ws1.Range("F5").Value = 295
$ws1.Range("F7").Value = 474
$ws1.Range("F10").Value = 303
$ws1.Range("F12").Value = 356
$ws1.Range("F13").Value = 5665
$ws1.Range("F14").Value = 73
$ws1.Range("F15").Value = 1749
$ws1.Range("F16").Value = 4104
$ws1.Range("F17").Value = 424
$ws1.Range("F20").Value = 4750
$ws1.Range("F21").Value = 6167
$ws1.Range("F23").Value = 1051
$ws1.Range("F24").Value = 683
$ws1.Range("F25").Value = 3743
$ws1.Range("F29").Value = 126
$ws1.Range("F31").Value = 1399
$ws1.Range("F33").Value = 541
$ws1.Range("F36").Value = 1698
$ws1.Range("F37").Value = 184
$ws1.Range("F39").Value = 1120
$ws1.Range("F40").Value = 34
$ws1.Range("F41").Value = 1342
$ws1.Range("F42").Value = 624
$ws1.Range("F45").Value = 127
$ws1.Range("F46").Value = 282

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 12
$ws2.Range("F24").Value = 71
$ws2.Range("F26").Value = 2
$ws2.Range("F27").Value = 50

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3823

# --- Sheet "全部类型" (All Types, aggregated) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3823
$ws4.Range("F3").Value = 2621
$ws4.Range("F4").Value = 568
$ws4.Range("F6").Value = 295
$ws4.Range("F10").Value = 474
$ws4.Range("F14").Value = 303
$ws4.Range("F16").Value = 356
$ws4.Range("F18").Value = 1749
$ws4.Range("F19").Value = 4751
$ws4.Range("F21").Value = 1051
$ws4.Range("F22").Value = 683
$ws4.Range("F23").Value = 3743
$ws4.Range("F27").Value = 126
$ws4.Range("F29").Value = 1399
$ws4.Range("F31").Value = 541
$ws4.Range("F35").Value = 1698
$ws4.Range("F37").Value = 1120
$ws4.Range("F39").Value = 624
$ws4.Range("F42").Value = 71
$ws4.Range("F45").Value = 127
$ws4.Range("F46").Value = 282
$ws4.Range("F48").Value = 2
